# Auto update Excel log
# Appends 7 new sensor-log rows (rows 31-37) to each of the three mmWave
# worksheets: mmWave(BR), mmWave(HR), mmWave(InBed).
#
# Columns: A=Date, B=Timestamp, C=Hour, D=Location, E=Value, F=Status

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$row,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        $value,
        [string]$status
    )

    # Column A holds a "YYYY-MM-DD" string that must stay literal text
    # (matching the existing rows in the sheet) rather than be coerced
    # into a date serial number, so force text formatting first.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

# ---------------------------------------------------------------------
# mmWave(InBed): Value column ("E") holds textual bed-occupancy state
# ---------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

Add-LogRow $wsInBed 31 "2026-01-28" "18:39:19" "18:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 32 "2026-01-28" "18:39:20" "18:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 33 "2026-01-28" "18:39:21" "18:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 34 "2026-01-28" "18:39:23" "18:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 35 "2026-01-28" "18:39:24" "18:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 36 "2026-01-28" "18:39:25" "18:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 37 "2026-01-28" "18:39:29" "18:00" "Bedroom" "In Bed" "Occupied"

# ---------------------------------------------------------------------
# mmWave(BR): Value column ("E") holds numeric breathing-rate readings
# ---------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

Add-LogRow $wsBR 31 "2026-01-28" "18:39:19" "18:00" "Bedroom" 0  "Occupied"
Add-LogRow $wsBR 32 "2026-01-28" "18:39:20" "18:00" "Bedroom" 26 "Occupied"
Add-LogRow $wsBR 33 "2026-01-28" "18:39:21" "18:00" "Bedroom" 2  "Occupied"
Add-LogRow $wsBR 34 "2026-01-28" "18:39:24" "18:00" "Bedroom" 53 "Occupied"
Add-LogRow $wsBR 35 "2026-01-28" "18:39:24" "18:00" "Bedroom" 33 "Occupied"
Add-LogRow $wsBR 36 "2026-01-28" "18:39:25" "18:00" "Bedroom" 2  "Occupied"
Add-LogRow $wsBR 37 "2026-01-28" "18:39:30" "18:00" "Bedroom" 1  "Occupied"

# ---------------------------------------------------------------------
# mmWave(HR): Value column ("E") holds numeric heart-rate readings
# ---------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

Add-LogRow $wsHR 31 "2026-01-28" "18:39:19" "18:00" "Bedroom" 0   "Occupied"
Add-LogRow $wsHR 32 "2026-01-28" "18:39:20" "18:00" "Bedroom" 74  "Occupied"
Add-LogRow $wsHR 33 "2026-01-28" "18:39:21" "18:00" "Bedroom" 50  "Occupied"
Add-LogRow $wsHR 34 "2026-01-28" "18:39:23" "18:00" "Bedroom" 101 "Occupied"
Add-LogRow $wsHR 35 "2026-01-28" "18:39:24" "18:00" "Bedroom" 81  "Occupied"
Add-LogRow $wsHR 36 "2026-01-28" "18:39:25" "18:00" "Bedroom" 50  "Occupied"
Add-LogRow $wsHR 37 "2026-01-28" "18:39:29" "18:00" "Bedroom" 49  "Occupied"
